# Trade #11 closed at 2026-02-17 07:53:33 - unknown UNKNOWN +0.000%
#
# Updates the "Summary" and "Strategy Status" sheets with the recalculated
# aggregate stats after closing trade #11, and appends the new trade row
# to both the "All Trades" and "MarketMaking" logs.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet: Total P&L %, Total Trades, Win Rate %
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B5").Value = -0.07000000000000001
$summary.Range("B6").Value = 11
$summary.Range("B9").Value = 36.36

# ---------------------------------------------------------------------
# Strategy Status sheet: MarketMaking row -> Trades, Win Rate %
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D4").Value = 11
$status.Range("G4").Value = 36.36

# ---------------------------------------------------------------------
# Append the new trade (#11) to "All Trades" and "MarketMaking" sheets.
# Row 11 is copied down to row 12 first so text-like columns (dates,
# strategy name, status, reasons, ...) keep their original text
# formatting, then only the cells that actually changed are overwritten.
# ---------------------------------------------------------------------
$sheetNames = @("All Trades", "MarketMaking")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("A11:Q11").Copy()
    $ws.Range("A12:Q12").PasteSpecial(-4104)

    $ws.Range("A12").Value = 11
    $ws.Range("C12").Value = "07:53:27"
    $ws.Range("F12").Value = 0.2
    $ws.Range("G12").Value = 0.2
    $ws.Range("I12").Value = 0
    $ws.Range("J12").Value = 0
    $ws.Range("Q12").Value = 0.11
}
